$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.596.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.53%  "

# Row 3
$ws.Range("D3").Value = "'1.830.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.40%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "'316.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

# Row 6
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
$ws.Range("E7").Value = "  +0.67%  "

# Row 8
$ws.Range("D8").Value = "'0.3989"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.87%  "

# Row 9
$ws.Range("D9").Value = "'0.07800"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.38%  "

# Row 10
$ws.Range("D10").Value = "'1.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "

# Row 11
$ws.Range("D11").Value = "'41.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.10%  "

# Row 12
$ws.Range("D12").Value = "'6.343"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.23%  "

# Row 13
$ws.Range("D13").Value = "'21.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "

# Row 14
$ws.Range("D14").Value = "'7.581"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.93%  "

# Row 15
$ws.Range("E15").Value = "  +0.11%  "

# Row 16
$ws.Range("D16").Value = "'1.832.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "

# Row 17
$ws.Range("D17").Value = "'93.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.21%  "

# Row 18
$ws.Range("D18").Value = "'0.00001091"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.42%  "

# Row 19
$ws.Range("D19").Value = "'0.06570"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").Value = "'17.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "

# Row 21
$ws.Range("E21").Value = "  +0.07%  "

# Row 22
$ws.Range("D22").Value = "'6.100"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.96%  "

# Row 23
$ws.Range("D23").Value = "'28.605.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.43%  "

# Row 24
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("E25").Value = "  +7.46%  "

# Row 26
$ws.Range("D26").Value = "'20.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'156.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "

# Row 28
$ws.Range("B28").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C28").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D28").Value = "'2.039.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

# Row 29
$ws.Range("D29").Value = "'2.424"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.89%  "

# Row 30
$ws.Range("D30").Value = "'125.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.85%  "

# Row 31
$ws.Range("D31").Value = "'1.142"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.23%  "

# Row 32
$ws.Range("D32").Value = "'0.1123"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "

# Row 33
$ws.Range("D33").Value = "'5.758"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.25%  "

# Row 34
$ws.Range("D34").Value = "'3.651"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "

# Row 35
$ws.Range("D35").Value = "'0.07294"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "

# Row 36
$ws.Range("D36").Value = "'0.2265"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "

# Row 37
$ws.Range("D37").Value = "'0.02348"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.33%  "

# Row 38
$ws.Range("D38").Value = "'8.946"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.24%  "

# Row 39
$ws.Range("E39").Value = "  +2.45%  "

# Row 40
$ws.Range("D40").Value = "'11.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.18%  "

# Row 41
$ws.Range("D41").Value = "'0.6309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "

# Row 42
$ws.Range("E42").Value = "  +1.36%  "

# Row 43
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("D44").Value = "'1.395"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.72%  "

# Row 45
$ws.Range("D45").Value = "'13.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.85%  "

# Row 46
$ws.Range("D46").Value = "'0.5925"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.69%  "

# Row 47
$ws.Range("D47").Value = "'3.714"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$ws.Range("D48").Value = "'125.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.11%  "

# Row 49
$ws.Range("E49").Value = "  +4.06%  "

# Row 50
$ws.Range("E50").Value = "  +0.95%  "

# Row 51
$ws.Range("D51").Value = "'0.06954"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.97%  "
